$d = $word.ActiveDocument

# 1) Merge "...que " + "este compra" + ", algo..." into a single run (removes the
#    gramStart/gramEnd proofErr markers Word had placed around "este compra").
$d.Content.Find.Execute("que este compra,", $false, $false, $false, $false, $false, $true, 1, $false, "que este compra,", 2) | Out-Null

# 2) Real text change: "un ejemplo cercano" -> "un objeto cercano"
$d.Content.Find.Execute("un ejemplo cercano", $false, $false, $false, $false, $false, $true, 1, $false, "un objeto cercano", 2) | Out-Null

# 3) Insert "del dataset " after "a cada punto "
$d.Content.Find.Execute("a cada punto dependiendo", $false, $false, $false, $false, $false, $true, 1, $false, "a cada punto del dataset dependiendo", 2) | Out-Null

# 4) Merge "normalizar los " + "atributos previo" + " a realizar " into one run
#    (removes the gramStart/gramEnd proofErr markers around "atributos previo").
$d.Content.Find.Execute("los atributos previo a realizar", $false, $false, $false, $false, $false, $true, 1, $false, "los atributos previo a realizar", 2) | Out-Null

# 5) Merge "Si, por el " + "contrario" + " se desea " into one run
#    (removes the gramStart/gramEnd proofErr markers around "contrario").
$d.Content.Find.Execute("Si, por el contrario se desea", $false, $false, $false, $false, $false, $true, 1, $false, "Si, por el contrario se desea", 2) | Out-Null

# 6) Merge "argumentar " + "que" + " para este " into one run
#    (removes the gramStart/gramEnd proofErr markers around "que").
$d.Content.Find.Execute("argumentar que para este", $false, $false, $false, $false, $false, $true, 1, $false, "argumentar que para este", 2) | Out-Null
